$wb = $excel.ActiveWorkbook

# --- Daily sheet updates (row 2) ---
$wsDaily = $wb.Worksheets.Item("Daily")
$wsDaily.Range("G2").Value = 2746.91
$wsDaily.Range("H2").Value = 6004.45
$wsDaily.Range("I2").Value = 700.27

# --- Hourly sheet updates (rows 9-19) ---
$wsHourly = $wb.Worksheets.Item("Hourly")
$wsHourly.Range("I9").Value = 30.81
$wsHourly.Range("H10").Value = 91.84999999999999
$wsHourly.Range("I10").Value = 397.17
$wsHourly.Range("J10").Value = 45.1
$wsHourly.Range("K10").Value = 22.96
$wsHourly.Range("M10").Value = 22.96
$wsHourly.Range("H11").Value = 229.5
$wsHourly.Range("I11").Value = 617.1
$wsHourly.Range("J11").Value = 69.59999999999999
$wsHourly.Range("H12").Value = 348.3
$wsHourly.Range("I12").Value = 719.67
$wsHourly.Range("J12").Value = 83.51000000000001
$wsHourly.Range("K12").Value = 87.06999999999999
$wsHourly.Range("M12").Value = 87.06999999999999
$wsHourly.Range("H13").Value = 428.19
$wsHourly.Range("I13").Value = 770.08
$wsHourly.Range("J13").Value = 91.14
$wsHourly.Range("K13").Value = 107.05
$wsHourly.Range("M13").Value = 107.05
$wsHourly.Range("H14").Value = 459.21
$wsHourly.Range("I14").Value = 787.04
$wsHourly.Range("J14").Value = 93.84
$wsHourly.Range("K14").Value = 114.8
$wsHourly.Range("M14").Value = 114.8
$wsHourly.Range("H15").Value = 437.83
$wsHourly.Range("I15").Value = 775.52
$wsHourly.Range("J15").Value = 91.98999999999999
$wsHourly.Range("H16").Value = 366.47
$wsHourly.Range("I16").Value = 732.1799999999999
$wsHourly.Range("J16").Value = 85.33
$wsHourly.Range("H17").Value = 253.82
$wsHourly.Range("I17").Value = 641.95
$wsHourly.Range("J17").Value = 72.77
$wsHourly.Range("K17").Value = 63.47
$wsHourly.Range("M17").Value = 63.47
$wsHourly.Range("H18").Value = 117.28
$wsHourly.Range("I18").Value = 454.47
$wsHourly.Range("J18").Value = 50.93
$wsHourly.Range("K18").Value = 29.35
$wsHourly.Range("M18").Value = 29.35
$wsHourly.Range("I19").Value = 78.45
